$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D: "success" flag.
# success = "1" when the list id (column A) is between 1 and 8 inclusive, else "0".

# Header cell D1: same text + same header style (border/bold/centered) as B1/C1.
$ws.Range("D1").Value = "success"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data cells D2:D25: plain text "0"/"1" values (no cell style).
for ($row = 2; $row -le 25; $row++) {
    $a = $ws.Cells.Item($row, 1).Value2
    if ($a -ge 1 -and $a -le 8) {
        $ws.Cells.Item($row, 4).Formula = "=""1"""
    } else {
        $ws.Cells.Item($row, 4).Formula = "=""0"""
    }
}

# Convert the formula-driven text results into plain static text values
# (paste-special values-only), matching how Excel stores literal text.
$dataRange = $ws.Range("D2:D25")
$dataRange.Copy()
$dataRange.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
